# Fruta / hortaliza, semanal
# Inserts 3 new rows at row 821 (shifting the existing 821-904 block down to
# 824-907) and populates the 3 new rows with the new weekly price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 821; Excel shifts rows 821:904 down to 824:907.
$ws.Rows("821:823").Insert()

# ---- Row 821: Naranja / Cara cara / Primera ----
$ws.Cells.Item(821, 1).Value = 10
$ws.Cells.Item(821, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(821, 3).Value = "La Araucanía"
$ws.Cells.Item(821, 4).Value = 44769
$ws.Cells.Item(821, 5).Value = 9
$ws.Cells.Item(821, 6).Value = "Fruta"
$ws.Cells.Item(821, 7).Value = 100102
$ws.Cells.Item(821, 8).Value = "Cítricos"
$ws.Cells.Item(821, 9).Value = 100102005
$ws.Cells.Item(821, 10).Value = "Naranja"
$ws.Cells.Item(821, 11).Value = "Cara cara"
$ws.Cells.Item(821, 12).Value = "Primera"
$ws.Cells.Item(821, 13).Value = 180
$ws.Cells.Item(821, 14).Value = 10000
$ws.Cells.Item(821, 15).Value = 10000
$ws.Cells.Item(821, 16).Value = 10000
$ws.Cells.Item(821, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(821, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(821, 19).Value = 667
$ws.Cells.Item(821, 20).Value = 15

# ---- Row 822: Naranja / Fukumoto / Primera ----
$ws.Cells.Item(822, 1).Value = 10
$ws.Cells.Item(822, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(822, 3).Value = "La Araucanía"
$ws.Cells.Item(822, 4).Value = 44769
$ws.Cells.Item(822, 5).Value = 9
$ws.Cells.Item(822, 6).Value = "Fruta"
$ws.Cells.Item(822, 7).Value = 100102
$ws.Cells.Item(822, 8).Value = "Cítricos"
$ws.Cells.Item(822, 9).Value = 100102005
$ws.Cells.Item(822, 10).Value = "Naranja"
$ws.Cells.Item(822, 11).Value = "Fukumoto"
$ws.Cells.Item(822, 12).Value = "Primera"
$ws.Cells.Item(822, 13).Value = 250
$ws.Cells.Item(822, 14).Value = 8000
$ws.Cells.Item(822, 15).Value = 8000
$ws.Cells.Item(822, 16).Value = 8000
$ws.Cells.Item(822, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(822, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(822, 19).Value = 533
$ws.Cells.Item(822, 20).Value = 15

# ---- Row 823: Naranja / Navel Late / Primera ----
$ws.Cells.Item(823, 1).Value = 10
$ws.Cells.Item(823, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(823, 3).Value = "La Araucanía"
$ws.Cells.Item(823, 4).Value = 44340
$ws.Cells.Item(823, 5).Value = 9
$ws.Cells.Item(823, 6).Value = "Fruta"
$ws.Cells.Item(823, 7).Value = 100102
$ws.Cells.Item(823, 8).Value = "Cítricos"
$ws.Cells.Item(823, 9).Value = 100102005
$ws.Cells.Item(823, 10).Value = "Naranja"
$ws.Cells.Item(823, 11).Value = "Navel Late"
$ws.Cells.Item(823, 12).Value = "Primera"
$ws.Cells.Item(823, 13).Value = 15
$ws.Cells.Item(823, 14).Value = 150000
$ws.Cells.Item(823, 15).Value = 150000
$ws.Cells.Item(823, 16).Value = 150000
$ws.Cells.Item(823, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(823, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(823, 19).Value = 375
$ws.Cells.Item(823, 20).Value = 400
